# Updates the cryptos list with the latest scraped data.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h); row 1 is the header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "66.319.80", "  -0.77%  "),
    @("Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.318.23", "  -0.67%  "),
    @("TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  +0.01%  "),
    @("BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "587.63", "  +2.49%  "),
    @("Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "183.52", "  +1.10%  "),
    @("XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.649", "  +8.50%  "),
    @("USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  +0.02%  "),
    @("LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "3.317.67", "  -0.68%  "),
    @("Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.126", "  -2.96%  "),
    @("Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "6.81", "  +2.23%  "),
    @("Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.401", "  -0.42%  "),
    @("WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.898.03", "  -0.58%  "),
    @("TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.132", "  -3.49%  "),
    @("WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "66.345.52", "  -0.78%  "),
    @("Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "26.16", "  -3.56%  "),
    @("WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.347.30", "  +0.41%  "),
    @("ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000164", "  -2.58%  "),
    @("BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "426.12", "  -2.58%  "),
    @("Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "5.53", "  -2.80%  "),
    @("Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "13.19", "  -3.04%  "),
    @("Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "7.40", "  -2.80%  "),
    @("Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "71.94", "  -2.63%  "),
    @("Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.02%  "),
    @("LEO", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo", "5.70", "  +0.37%  "),
    @("WrappedeETH", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth", "3.464.44", "  -0.66%  "),
    @("Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.515", "  -0.84%  "),
    @("Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.204", "  +6.93%  "),
    @("PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0000114", "  -3.96%  "),
    @("InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "8.96", "  -0.89%  "),
    @("Binance-PegBSC-USD", "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd", "0.999", "  -0.04%  "),
    @("PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.93", "  -1.45%  "),
    @("EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "22.34", "  -2.21%  "),
    @("USDe", "https://coinranking.com/coin/exbfr2U-0+usde-usde", "1.00", "  +0.09%  "),
    @("NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "5.18", "  -2.45%  "),
    @("Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.59", "  -3.06%  "),
    @("Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "1.18", "  -4.09%  "),
    @("Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "159.92", "  -1.28%  "),
    @("ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.44", "  -3.77%  "),
    @("Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "2.885.39", "  +1.90%  "),
    @("Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "1.80", "  -2.27%  "),
    @("EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "26.45", "  -5.49%  "),
    @("Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.764", "  -3.44%  "),
    @("Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.32", "  -2.80%  "),
    @("OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "40.11", "  -0.03%  "),
    @("Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0664", "  -1.05%  "),
    @("RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "5.93", "  -5.79%  "),
    @("dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "2.31", "  -2.32%  "),
    @("InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "23.27", "  -5.21%  "),
    @("Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "313.43", "  -3.06%  "),
    @("VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0272", "  -0.42%  ")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $row++
}
